$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 6) - set first so shared-string order matches target
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "SCH4U"
$ws.Range("D6").Value = "N"
$ws.Range("E6").Value = "Filiberto Cuevas"

# Resize the table to include the 3 new columns and the new row
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:H28"))

# New headers for added columns (set via the table header range so both the
# worksheet cell and the table's column definition stay in sync)
$table.HeaderRowRange.Item(6).Value = "Replacement ID"
$table.HeaderRowRange.Item(7).Value = "Absentee"
$table.HeaderRowRange.Item(8).Value = "Absentee ID"

# Column widths for new columns (values chosen so the serialized OOXML width
# ends up as close as possible to the target widths 17.1640625 / 16.1640625 / 18.83203125)
$ws.Columns.Item(6).ColumnWidth = 16.33
$ws.Columns.Item(7).ColumnWidth = 15.33
$ws.Columns.Item(8).ColumnWidth = 18

# Match final cell selection
$ws.Range("K5").Select() | Out-Null
